$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text without altering its number format/style,
# even when the text looks like a number (e.g. "39.49", "0.08020").
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "20.549.03"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.472.60"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  -0.25%  "
Set-TextValue $ws.Range("D5") "0.9595"
$ws.Range("E5").Value = "  +5.85%  "
Set-TextValue $ws.Range("D6") "277.37"
$ws.Range("E6").Value = "  +0.01%  "
Set-TextValue $ws.Range("D7") "0.3618"
$ws.Range("E7").Value = "  -0.80%  "
Set-TextValue $ws.Range("D8") "0.3075"
$ws.Range("E8").Value = "  -1.14%  "
Set-TextValue $ws.Range("D9") "39.49"
$ws.Range("E9").Value = "  +1.19%  "
Set-TextValue $ws.Range("D10") "1.077"
$ws.Range("E10").Value = "  +5.94%  "
Set-TextValue $ws.Range("D11") "0.06652"
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("E12").Value = "  -0.09%  "
Set-TextValue $ws.Range("D13") "5.517"
$ws.Range("E13").Value = "  +2.62%  "
Set-TextValue $ws.Range("D14") "18.16"
$ws.Range("E14").Value = "  +3.44%  "
Set-TextValue $ws.Range("D15") "6.170"
$ws.Range("E15").Value = "  +2.10%  "
Set-TextValue $ws.Range("D16") "0.9594"
$ws.Range("E16").Value = "  +1.81%  "
Set-TextValue $ws.Range("D17") "0.00001026"
$ws.Range("D18").Value = "1.473.24"
$ws.Range("E18").Value = "  +2.28%  "
Set-TextValue $ws.Range("D19") "0.05934"
$ws.Range("E19").Value = "  +5.36%  "
$ws.Range("E20").Value = "  +1.76%  "
Set-TextValue $ws.Range("D21") "5.505"
$ws.Range("E21").Value = "  +2.22%  "
Set-TextValue $ws.Range("D22") "14.57"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("E23").Value = "  +3.72%  "
Set-TextValue $ws.Range("D24") "2.263"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "20.554.49"
$ws.Range("E25").Value = "  +1.43%  "
Set-TextValue $ws.Range("D26") "143.33"
$ws.Range("E26").Value = "  +3.82%  "
Set-TextValue $ws.Range("D27") "2.127"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "1.633.16"
$ws.Range("E29").Value = "  +2.74%  "
Set-TextValue $ws.Range("D30") "113.97"
$ws.Range("E30").Value = "  +3.91%  "
Set-TextValue $ws.Range("D31") "3.891"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "4.957"
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D33") "0.08020"
$ws.Range("E33").Value = "  +4.40%  "
Set-TextValue $ws.Range("D34") "0.8057"
$ws.Range("E34").Value = "  +0.67%  "
Set-TextValue $ws.Range("D35") "1.514"
$ws.Range("E35").Value = "  +4.80%  "
Set-TextValue $ws.Range("D36") "1.214"
$ws.Range("E36").Value = "  +6.89%  "
Set-TextValue $ws.Range("D37") "0.05775"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("E39").Value = "  +3.83%  "
Set-TextValue $ws.Range("D40") "10.42"
$ws.Range("E40").Value = "  +2.89%  "
Set-TextValue $ws.Range("D41") "0.9596"
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("E42").Value = "  +2.14%  "
Set-TextValue $ws.Range("D43") "7.430"
$ws.Range("E43").Value = "  +4.74%  "
Set-TextValue $ws.Range("D44") "0.5281"
$ws.Range("E44").Value = "  +1.20%  "
Set-TextValue $ws.Range("D45") "3.521"
$ws.Range("E45").Value = "  +0.21%  "
Set-TextValue $ws.Range("D46") "12.19"
$ws.Range("E46").Value = "  +1.89%  "
Set-TextValue $ws.Range("D47") "118.71"
$ws.Range("E47").Value = "  -0.41%  "
Set-TextValue $ws.Range("D48") "0.5204"
$ws.Range("E48").Value = "  +1.63%  "
Set-TextValue $ws.Range("D49") "1.821"
$ws.Range("E49").Value = "  +3.82%  "
Set-TextValue $ws.Range("D50") "0.06468"
$ws.Range("E50").Value = "  +2.30%  "
Set-TextValue $ws.Range("D51") "0.9840"
$ws.Range("E51").Value = "  -0.73%  "
